$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r2 = $ws.Range("D2:E2")
$r2.NumberFormat = "@"
$ws.Range("D2").Value = '254.84'
$ws.Range("E2").Value = '3.66%'
$r2.Style = "Normal"

$r3 = $ws.Range("E3")
$r3.NumberFormat = "@"
$ws.Range("E3").Value = '-6.95%'
$r3.Style = "Normal"

$r4 = $ws.Range("D4:E4")
$r4.NumberFormat = "@"
$ws.Range("D4").Value = '5.241'
$ws.Range("E4").Value = '1.69%'
$r4.Style = "Normal"

$r5 = $ws.Range("D5:E5")
$r5.NumberFormat = "@"
$ws.Range("D5").Value = '0.05884'
$ws.Range("E5").Value = '2.09%'
$r5.Style = "Normal"

$r6 = $ws.Range("D6:E6")
$r6.NumberFormat = "@"
$ws.Range("D6").Value = '6.736'
$ws.Range("E6").Value = '1.04%'
$r6.Style = "Normal"

$r7 = $ws.Range("D7:E7")
$r7.NumberFormat = "@"
$ws.Range("D7").Value = '0.8655'
$ws.Range("E7").Value = '1.65%'
$r7.Style = "Normal"

$r8 = $ws.Range("D8:E8")
$r8.NumberFormat = "@"
$ws.Range("D8").Value = '0.9885'
$ws.Range("E8").Value = '15.31%'
$r8.Style = "Normal"

$r9 = $ws.Range("D9:E9")
$r9.NumberFormat = "@"
$ws.Range("D9").Value = '0.1409'
$ws.Range("E9").Value = '2.05%'
$r9.Style = "Normal"

$r10 = $ws.Range("D10:E10")
$r10.NumberFormat = "@"
$ws.Range("D10").Value = '0.07163'
$ws.Range("E10").Value = '1.13%'
$r10.Style = "Normal"

$r11 = $ws.Range("D11:E11")
$r11.NumberFormat = "@"
$ws.Range("D11").Value = '0.03187'
$ws.Range("E11").Value = '-2.30%'
$r11.Style = "Normal"

$r12 = $ws.Range("D12:E12")
$r12.NumberFormat = "@"
$ws.Range("D12").Value = '0.09225'
$ws.Range("E12").Value = '-1.48%'
$r12.Style = "Normal"

$r13 = $ws.Range("D13:E13")
$r13.NumberFormat = "@"
$ws.Range("D13").Value = '0.001545'
$ws.Range("E13").Value = '0.43%'
$r13.Style = "Normal"

$r14 = $ws.Range("D14:E14")
$r14.NumberFormat = "@"
$ws.Range("D14").Value = '0.01055'
$ws.Range("E14").Value = '1,667.39%'
$r14.Style = "Normal"

$r15 = $ws.Range("D15:E15")
$r15.NumberFormat = "@"
$ws.Range("D15").Value = '0.005838'
$ws.Range("E15").Value = '-5.03%'
$r15.Style = "Normal"

$r16 = $ws.Range("E16")
$r16.NumberFormat = "@"
$ws.Range("E16").Value = '-0.78%'
$r16.Style = "Normal"

$r17 = $ws.Range("D17:E17")
$r17.NumberFormat = "@"
$ws.Range("D17").Value = '3.223'
$ws.Range("E17").Value = '-0.83%'
$r17.Style = "Normal"

$r18 = $ws.Range("E18")
$r18.NumberFormat = "@"
$ws.Range("E18").Value = '-0.64%'
$r18.Style = "Normal"

$r19 = $ws.Range("E19")
$r19.NumberFormat = "@"
$ws.Range("E19").Value = '1.74%'
$r19.Style = "Normal"

$r20 = $ws.Range("D20:E20")
$r20.NumberFormat = "@"
$ws.Range("D20").Value = '0.03473'
$ws.Range("E20").Value = '2.34%'
$r20.Style = "Normal"

$r21 = $ws.Range("D21:E21")
$r21.NumberFormat = "@"
$ws.Range("D21").Value = '0.1320'
$ws.Range("E21").Value = '-0.26%'
$r21.Style = "Normal"

$r22 = $ws.Range("E22")
$r22.NumberFormat = "@"
$ws.Range("E22").Value = '1.94%'
$r22.Style = "Normal"

$r23 = $ws.Range("D23:E23")
$r23.NumberFormat = "@"
$ws.Range("D23").Value = '0.04159'
$ws.Range("E23").Value = '0.97%'
$r23.Style = "Normal"

$r24 = $ws.Range("E24")
$r24.NumberFormat = "@"
$ws.Range("E24").Value = '-2.04%'
$r24.Style = "Normal"

$r25 = $ws.Range("D25:E25")
$r25.NumberFormat = "@"
$ws.Range("D25").Value = '0.001224'
$ws.Range("E25").Value = '-0.03%'
$r25.Style = "Normal"

$r26 = $ws.Range("D26:E26")
$r26.NumberFormat = "@"
$ws.Range("D26").Value = '0.004795'
$ws.Range("E26").Value = '15.52%'
$r26.Style = "Normal"

$r27 = $ws.Range("D27:E27")
$r27.NumberFormat = "@"
$ws.Range("D27").Value = '0.0001200'
$ws.Range("E27").Value = '0.08%'
$r27.Style = "Normal"

$r28 = $ws.Range("E28")
$r28.NumberFormat = "@"
$ws.Range("E28").Value = '1.26%'
$r28.Style = "Normal"

$r40 = $ws.Range("D40:E40")
$r40.NumberFormat = "@"
$ws.Range("D40").Value = '0.03811'
$ws.Range("E40").Value = '1.38%'
$r40.Style = "Normal"

$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$r41 = $ws.Range("D41:E41")
$r41.NumberFormat = "@"
$ws.Range("D41").Value = '0.005732'
$ws.Range("E41").Value = '-0.10%'
$r41.Style = "Normal"

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$r42 = $ws.Range("D42:E42")
$r42.NumberFormat = "@"
$ws.Range("D42").Value = '0.1100'
$ws.Range("E42").Value = '2.80%'
$r42.Style = "Normal"

$r43 = $ws.Range("D43:E43")
$r43.NumberFormat = "@"
$ws.Range("D43").Value = '0.002343'
$ws.Range("E43").Value = '6.57%'
$r43.Style = "Normal"

$r44 = $ws.Range("D44:E44")
$r44.NumberFormat = "@"
$ws.Range("D44").Value = '0.009697'
$ws.Range("E44").Value = '8.72%'
$r44.Style = "Normal"

$r45 = $ws.Range("D45:E45")
$r45.NumberFormat = "@"
$ws.Range("D45").Value = '0.00005239'
$ws.Range("E45").Value = '-4.31%'
$r45.Style = "Normal"

$r46 = $ws.Range("E46")
$r46.NumberFormat = "@"
$ws.Range("E46").Value = '0.08%'
$r46.Style = "Normal"

$r47 = $ws.Range("D47:E47")
$r47.NumberFormat = "@"
$ws.Range("D47").Value = '0.09303'
$ws.Range("E47").Value = '31.09%'
$r47.Style = "Normal"

$r48 = $ws.Range("D48:E48")
$r48.NumberFormat = "@"
$ws.Range("D48").Value = '0.002146'
$ws.Range("E48").Value = '-12.99%'
$r48.Style = "Normal"

$r49 = $ws.Range("D49:E49")
$r49.NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'
$ws.Range("E49").Value = '0.08%'
$r49.Style = "Normal"

$r50 = $ws.Range("D50:E50")
$r50.NumberFormat = "@"
$ws.Range("D50").Value = '0.0002001'
$ws.Range("E50").Value = '0.08%'
$r50.Style = "Normal"
